$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert two new rows right after row 4 (date field), to put the date on its own screen.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 2).Value = "end screen"
$ws.Cells.Item(6, 2).Value = "begin screen"

# Activate the survey sheet/tab and select B12 (the new "begin screen" row for the
# reporting-year/month/module screen) to match the saved selection state.
$ws.Activate()
$ws.Range("B12").Select()
